$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.753.06'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '2.276.40'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '251.01'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').Value = '0.644'
$ws.Range('E6').Value = '  +2.89%  '
$ws.Range('D7').Value = '75.94'
$ws.Range('E7').Value = '  +6.66%  '
$ws.Range('D9').Value = '0.644'
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('D10').Value = '40.14'
$ws.Range('E10').Value = '  +3.04%  '
$ws.Range('D11').Value = '0.0975'
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('D12').Value = '7.37'
$ws.Range('E12').Value = '  -1.35%  '
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '2.620.62'
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('D15').Value = '15.01'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('D16').Value = '0.865'
$ws.Range('E16').Value = '  -2.19%  '
$ws.Range('D17').Value = '2.286.87'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').Value = '42.675.80'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').Value = '0.0₃0995'
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').Value = '6.20'
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('D21').Value = '72.28'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D22').Value = '236.38'
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('D23').Value = '2.15'
$ws.Range('E23').Value = '  +3.52%  '
$ws.Range('E24').Value = '  -2.21%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').Value = '11.23'
$ws.Range('E26').Value = '  -2.18%  '
$ws.Range('D27').Value = '2.39'
$ws.Range('E27').Value = '  -1.79%  '
$ws.Range('D28').Value = '2.14'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('D29').Value = '167.52'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').Value = '20.96'
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '6.43'
$ws.Range('E31').Value = '  -2.04%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0856'
$ws.Range('E32').Value = '  +7.54%  '
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('D34').Value = '31.89'
$ws.Range('E35').Value = '  +1.64%  '
$ws.Range('E36').Value = '  +2.30%  '
$ws.Range('D37').Value = '4.75'
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('E38').Value = '  -4.74%  '
$ws.Range('D39').Value = '13.51'
$ws.Range('E39').Value = '  +8.13%  '
$ws.Range('D40').Value = '2.28'
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('D41').Value = '5.88'
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('E42').Value = '  +2.48%  '
$ws.Range('D43').Value = '61.50'
$ws.Range('E43').Value = '  -1.26%  '
$ws.Range('D44').Value = '8.91'
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('D45').Value = '106.25'
$ws.Range('E45').Value = '  +11.83%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').Value = '4.71'
$ws.Range('E46').Value = '  -2.50%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.100'
$ws.Range('E47').Value = '  -1.67%  '
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('E50').Value = '  -2.39%  '
$ws.Range('D51').Value = '4.19'
$ws.Range('E51').Value = '  -2.35%  '
